$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tableau de Bord")
$ws.Rows("12:13").Insert()
$wb.Names.Add("BASE_URL", "='Tableau de Bord'!$D$14")
$wb.Names.Add("VERSION_MSTG", "='Tableau de Bord'!$D$13")
foreach ($n in $wb.Names) {
    Write-Host ("Name: " + $n.Name + " => " + $n.RefersTo)
}
